$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.169968962669373
$ws.Range("B1").Value = 2.379768371582031
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.37913966178894
$ws.Range("E1").Value = 1.210723996162415
